# household.xlsx regen: update form_version setting and refresh the
# active-sheet/selection UI state left behind by the previous save.

$wb = $excel.ActiveWorkbook

# --- settings sheet: form_version value bumped from 1 to 20130408 ---
$settings = $wb.Worksheets.Item("settings")
$settings.Range("B3").Value = 20130408

# --- restore selections on each sheet, then make "settings" the active tab ---
$survey = $wb.Worksheets.Item("survey")
$survey.Range("C20").Select()

$choices = $wb.Worksheets.Item("choices")
$choices.Range("E12").Select()

$settings.Activate()
$settings.Range("B7").Select()
